# Update the "StdDev_*" weighted-mean (WMean) raw values in row 3.
# These are the only literal (non-formula) values in the diff; all the
# dependent "HEREstats" rows (6-11) reference these cells through
# formulas (e.g. C6=G3, C7=M3, ...) and will recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value  = 11.0147556915277
$ws.Range("M3").Value  = 10.845823445326401
$ws.Range("R3").Value  = 11.7141004374893
$ws.Range("Y3").Value  = 10.6532579483576
$ws.Range("AE3").Value = 10.6817822447937
$ws.Range("AK3").Value = 10.09293556908
$ws.Range("AQ3").Value = 10.6931989366949

# The sheet view had scrolled so column I was the left-most visible
# column (topLeftCell="I1"); restore the default scroll position so the
# view starts at A1 again (this also drops the now-stale topLeftCell
# attribute from the saved XML).
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 1
